$wb = $excel.ActiveWorkbook

# --- FUSELAGE sheet: swap the TORENBEEK_1982 / SFORZA comparison rows (23 <-> 24) ---
$wsFuselage = $wb.Worksheets.Item("FUSELAGE")
$wsFuselage.Range("A23").Value = "SFORZA"
$wsFuselage.Range("C23").Value = 17.143322222222217
$wsFuselage.Range("A24").Value = "TORENBEEK_1982"
$wsFuselage.Range("C24").Value = 16.8345

# --- WING sheet: swap both the Xcg (23<->24) and Ycg (27<->28) comparison rows ---
$wsWing = $wb.Worksheets.Item("WING")
$wsWing.Range("A23").Value = "SFORZA"
$wsWing.Range("C23").Value = 4.3631082000119275
$wsWing.Range("A24").Value = "TORENBEEK_1982"
$wsWing.Range("C24").Value = 3.5939754358446514

$wsWing.Range("A27").Value = "SFORZA"
$wsWing.Range("C27").Value = 4.998846772296348
$wsWing.Range("A28").Value = "TORENBEEK_1982"
$wsWing.Range("C28").Value = 6.114221148470394
